$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (F column) counters
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5438
$ws1.Range("F6").Value = 78

# Sheet "全部类型": same rows duplicated, update corresponding cells
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5438
$ws4.Range("F7").Value = 78
